$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "27.546.25"
$ws.Range("E2").Value = "  +1.95%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.864.71"
$ws.Range("E3").Value = "  +0.96%  "

# Row 4 (TetherUSD)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.53"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6 (USDC)
$ws.Range("E6").Value = "  -0.22%  "

# Row 7 (XRP)
$ws.Range("E7").Value = "  +0.50%  "

# Row 8 (Cardano)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3805"
$ws.Range("E8").Value = "  +3.46%  "

# Row 9 (Dogecoin)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07345"
$ws.Range("E9").Value = "  +1.43%  "

# Row 10 (Polygon)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9342"
$ws.Range("E10").Value = "  +0.45%  "

# Row 11 (Solana)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.77"
$ws.Range("E11").Value = "  +4.68%  "

# Row 12 (TRON)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07788"

# Row 13 (WrappedEther)
$ws.Range("D13").Value = "1.897.79"
$ws.Range("E13").Value = "  +3.12%  "

# Row 14 (Polkadot)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.443"
$ws.Range("E14").Value = "  +0.95%  "

# Row 15 (Chainlink)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.576"
$ws.Range("E15").Value = "  +1.60%  "

# Row 16 (Litecoin)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.36"
$ws.Range("E16").Value = "  +1.60%  "

# Row 17 (BinanceUSD)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.014"
$ws.Range("E17").Value = "  -0.21%  "

# Row 18 (ShibaInu)
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008831"
$ws.Range("E18").Value = "  +1.89%  "

# Row 19
$ws.Range("E19").Value = "  -0.26%  "

# Row 20
$ws.Range("D20").Value = "27.548.79"
$ws.Range("E20").Value = "  +1.86%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.70"
$ws.Range("E21").Value = "  +1.02%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.104"
$ws.Range("E22").Value = "  +1.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.73"
$ws.Range("E23").Value = "  +0.71%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.936"
$ws.Range("E24").Value = "  +0.40%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.21"
$ws.Range("E25").Value = "  +2.20%  "

# Row 26
$ws.Range("E26").Value = "  +1.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.028"
$ws.Range("E27").Value = "  +2.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.45"
$ws.Range("E28").Value = "  +0.84%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.959"
$ws.Range("E29").Value = "  +0.18%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08883"
$ws.Range("E30").Value = "  +0.20%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.331"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.208"
$ws.Range("E32").Value = "  +2.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7565"
$ws.Range("E33").Value = "  +2.68%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.596"
$ws.Range("E34").Value = "  +1.93%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.703"
$ws.Range("E35").Value = "  +1.89%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.124"
$ws.Range("E36").Value = "  +0.52%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02041"
$ws.Range("E37").Value = "  +3.54%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5676"
$ws.Range("E38").Value = "  +7.88%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05356"
$ws.Range("E39").Value = "  +1.85%  "

# Row 40
$ws.Range("E40").Value = "  +0.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.029"
$ws.Range("E41").Value = "  -0.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.526"
$ws.Range("E42").Value = "  +2.95%  "

# Row 43
$ws.Range("E43").Value = "  +0.66%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.73"
$ws.Range("E44").Value = "  +1.83%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4886"

# Row 46
$ws.Range("E46").Value = "  -0.31%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.01"
$ws.Range("E47").Value = "  +3.33%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.667"
$ws.Range("E48").Value = "  +3.64%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.45"
$ws.Range("E49").Value = "  +2.67%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06091"
$ws.Range("E50").Value = "  +0.48%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9106"
$ws.Range("E51").Value = "  +2.06%  "
